# Remove the review row for armonravid@gmail.com / avishaybar12@gmail.com
# ("I MADE IT TO THE LAST LEVEL AND FAILED...") which lived in row 2.
# Deleting the row shifts every subsequent row up by one and drops the
# now-orphaned hyperlink that was anchored to the old C2 cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Hyperlinks.Delete()
$ws.Rows.Item(2).Delete()
